$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 71's formatting: it is no longer the last row, so its border
# style switches from the "last row" variant (no border) to the "interior
# row" variant (thin top+bottom), matching rows like 60/64/68. We replicate
# this by copying formats from row 68 (style 9/9/10/10/10) onto row 71. ---
$ws.Range("A68:E68").Copy()
$ws.Range("A71:E71").PasteSpecial(-4122)

# --- Build new row 72. Its target style mix (A/B=9, C=10, D/E=8) does not
# exist verbatim on any single existing row, so stitch it from two sources:
# A:C from row 68 (9,9,10) and D:E from row 67 (8,8). ---
$ws.Range("A68:C68").Copy()
$ws.Range("A72:C72").PasteSpecial(-4122)
$ws.Range("D67:E67").Copy()
$ws.Range("D72:E72").PasteSpecial(-4122)

$ws.Range("A72").Value = 'SCRIPT/G01P03A/us2206.ssb'
$ws.Range("B72").Value = 18
$ws.Range("C72").Value = ' Hey, hey! Let\''s have another fun\nday of exploring, hey, hey!'
$ws.Range("D72").Value = ' Эй, эй! Самое время для ещё одного\nувлекательного дня исследований, эй, эй!'
$ws.Range("E72").Value = ' Üê, üê! Òàíïå âñåíÿ äìÿ åþæ ïäîïãï\nôâìåëàóåìûîïãï äîÿ éòòìåäïâàîéê, üê, üê!'

# Row 72 mirrors the wrapped two-line height used throughout this sheet.
$ws.Rows.Item(72).RowHeight = 43.2

# Scroll position / active selection, per the saved view state.
$ws.Range("D71").Select()

Write-Host "Applied Korfish row-72 edit"
